$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Average GPU utilization" header from G1 to F1 ---
$headerText = $ws.Range("G1").Value2
$ws.Range("G1").ClearContents()
$ws.Range("F1").Value = $headerText

# --- Resize the columns around the relocated / new tables ---
# (Excel/COM quantizes ColumnWidth to whole-pixel increments; these are the
# input values that land closest to the authored 20.92 / 23.89 character
# widths once re-expanded through the pixel->width formula.)
$ws.Columns("F").ColumnWidth = 20.17
$ws.Columns("G").ColumnWidth = 23.0

# --- Highlight "low recall" (C23) in red to flag the mixed precision issue ---
$ws.Range("C23").Font.Color = 1974729

# --- Add the new "1024, mixed precision, XLA" mirrored table next to the
#     existing "batch_size = 1024" table (rows 30-35, columns G-J) ---
$ws.Range("G30").Value = "1024, mixed precision, XLA"

$ws.Range("H31").Value = 64
$ws.Range("I31").Value = 128
$ws.Range("J31").Value = 256

$ws.Range("G32").Value = 4
$ws.Range("I32").Value = 123.31

$ws.Range("G33").Value = 3

$ws.Range("G34").Value = 2
$ws.Range("I34").Value = 112.78

$ws.Range("G35").Value = 1

# --- Update the view: scroll + selection to the new table ---
$ws.Range("I24").Select() | Out-Null
